# Scheduled market-data refresh: update current average prices / leve profit
# figures across the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18969.084
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 18969.084
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 18969.084
$ws.Range("N32").Value = -19543.084
$ws.Range("M32").ClearContents()

$ws.Range("H110").Value = 896971.5600000001
$ws.Range("I110").Value = 896971.5600000001
$ws.Range("K110").Value = 896971.5600000001
$ws.Range("M110").Value = -894926.5600000001

$ws.Range("H122").Value = 10419150
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 10419150
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 31257450
$ws.Range("N122").Value = -31262350
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 25053752
$ws.Range("J86").Value = 100000
$ws.Range("L86").Value = 100000
$ws.Range("N86").Value = -102246

$ws.Range("H89").Value = 25053752
$ws.Range("J89").Value = 100000
$ws.Range("L89").Value = 500000
$ws.Range("N89").Value = -511232

$ws.Range("H94").Value = 2781820.8
$ws.Range("I94").Value = 3705305
$ws.Range("K94").Value = 3705305
$ws.Range("M94").Value = -3704854

$ws.Range("H107").Value = 2978073
$ws.Range("I107").Value = 3760972.5
$ws.Range("J107").Value = 3054.6
$ws.Range("K107").Value = 3760972.5
$ws.Range("L107").Value = 3054.6
$ws.Range("M107").Value = -3759052.5
$ws.Range("N107").Value = -6894.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6000
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 6500
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 6500
$ws.Range("M99").Value = -4002
$ws.Range("N99").Value = -9496

$ws.Range("H107").Value = 3413.318
$ws.Range("I107").Value = 2721.4546
$ws.Range("K107").Value = 2721.4546
$ws.Range("M107").Value = -801.4546

$ws.Range("H122").Value = 1788.4359
$ws.Range("I122").Value = 1677.9
$ws.Range("J122").Value = 2156.889
$ws.Range("K122").Value = 5033.700000000001
$ws.Range("L122").Value = 6470.667
$ws.Range("M122").Value = -2583.700000000001
$ws.Range("N122").Value = -11370.667

$ws.Range("H125").Value = 39165
$ws.Range("J125").Value = 39165
$ws.Range("L125").Value = 39165
$ws.Range("N125").Value = -44085

$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -24440

$ws.Range("H132").Value = 88813
$ws.Range("I132").Value = 93545.27
$ws.Range("J132").Value = 84080.73
$ws.Range("K132").Value = 280635.81
$ws.Range("L132").Value = 252242.19
$ws.Range("M132").Value = -278105.81
$ws.Range("N132").Value = -257302.19

$ws.Range("H134").Value = 33326.7
$ws.Range("I134").Value = 62386
$ws.Range("J134").Value = 4267.4
$ws.Range("K134").Value = 187158
$ws.Range("L134").Value = 12802.2
$ws.Range("M134").Value = -184623
$ws.Range("N134").Value = -17872.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 262.5
$ws.Range("I86").Value = 100
$ws.Range("J86").Value = 425
$ws.Range("K86").Value = 300
$ws.Range("L86").Value = 1275
$ws.Range("M86").Value = 886
$ws.Range("N86").Value = -3647

$ws.Range("H89").Value = 262.5
$ws.Range("I89").Value = 100
$ws.Range("J89").Value = 425
$ws.Range("K89").Value = 900
$ws.Range("L89").Value = 3825
$ws.Range("M89").Value = 5028
$ws.Range("N89").Value = -15681

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7160
$ws.Range("I55").Value = 3600
$ws.Range("K55").Value = 3600
$ws.Range("M55").Value = -3273

$ws.Range("H64").Value = 29999
$ws.Range("J64").Value = 29999
$ws.Range("L64").Value = 29999
$ws.Range("N64").Value = -30495

$ws.Range("H67").Value = 29999
$ws.Range("J67").Value = 29999
$ws.Range("L67").Value = 29999
$ws.Range("N67").Value = -31715

$ws.Range("H107").Value = 628.2857
$ws.Range("I107").Value = 1362.5
$ws.Range("J107").Value = 334.6
$ws.Range("K107").Value = 1362.5
$ws.Range("L107").Value = 334.6
$ws.Range("M107").Value = 557.5
$ws.Range("N107").Value = -4174.6

$ws.Range("H131").Value = 71975
$ws.Range("J131").Value = 71975
$ws.Range("L131").Value = 71975
$ws.Range("N131").Value = -82055

$ws.Range("H132").Value = 3326.4443
$ws.Range("I132").Value = 3251.8096
$ws.Range("J132").Value = 3587.6667
$ws.Range("K132").Value = 9755.4288
$ws.Range("L132").Value = 10763.0001
$ws.Range("M132").Value = -7225.4288
$ws.Range("N132").Value = -15823.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3537.7
$ws.Range("I40").Value = 2609.12
$ws.Range("J40").Value = 8180.6
$ws.Range("K40").Value = 2609.12
$ws.Range("L40").Value = 8180.6
$ws.Range("M40").Value = -2473.12
$ws.Range("N40").Value = -8452.6

$ws.Range("H43").Value = 6266.6665
$ws.Range("I43").Value = 6266.6665
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 6266.6665
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -6073.6665
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 4354.524
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 7358.636
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 7358.636
$ws.Range("M46").Value = -862
$ws.Range("N46").Value = -7734.636

$ws.Range("H55").Value = 2070.0588
$ws.Range("I55").Value = 2005
$ws.Range("J55").Value = 2143.25
$ws.Range("K55").Value = 2005
$ws.Range("L55").Value = 2143.25
$ws.Range("M55").Value = -1832
$ws.Range("N55").Value = -2489.25

$ws.Range("H61").Value = 10112121
$ws.Range("I61").Value = 15883263
$ws.Range("J61").Value = 12623.75
$ws.Range("K61").Value = 15883263
$ws.Range("L61").Value = 12623.75
$ws.Range("M61").Value = -15883061
$ws.Range("N61").Value = -13027.75

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 10112121
$ws.Range("I113").Value = 15883263
$ws.Range("J113").Value = 12623.75
$ws.Range("K113").Value = 15883263
$ws.Range("L113").Value = 12623.75
$ws.Range("M113").Value = -15881093
$ws.Range("N113").Value = -16963.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13864
$ws.Range("J74").Value = 13864
$ws.Range("L74").Value = 13864
$ws.Range("N74").Value = -15736

$ws.Range("H77").Value = 13864
$ws.Range("J77").Value = 13864
$ws.Range("L77").Value = 41592
$ws.Range("N77").Value = -50952

$ws.Range("H122").Value = 5882.6
$ws.Range("J122").Value = 6572.8
$ws.Range("L122").Value = 19718.4
$ws.Range("N122").Value = -24618.4

$ws.Range("H126").Value = 2256.1
$ws.Range("I126").Value = 2401.7778
$ws.Range("K126").Value = 7205.3334
$ws.Range("M126").Value = -4735.3334

$ws.Range("H129").Value = 58900
$ws.Range("J129").Value = 73733.336
$ws.Range("L129").Value = 73733.336
$ws.Range("N129").Value = -83733.336

$ws.Range("H131").Value = 105000
$ws.Range("J131").Value = 105000
$ws.Range("L131").Value = 105000
$ws.Range("N131").Value = -115080
